$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: split the old "Nama RotiPRM" header into two
# separate headers ("Nama Roti" and "PRM"), shifting PRS/PRD over.
$ws.Range("C1").Value = "PRM"
$ws.Range("B1").Value = "Nama Roti"
$ws.Range("D1").Value = "PRS"
$ws.Range("E1").Value = "PRD"

# Add the new (erroring) formula in C27
$ws.Range("C27").Formula = "=av"

# Update the selection to match the edited workbook
$ws.Range("C1").Select()
